$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing columns A,B shift to B,C
$ws.Columns.Item(1).Insert() | Out-Null

# Set the new column A width (closest achievable quantized value to 35.85546875)
$ws.Columns.Item(1).ColumnWidth = 35

# Header for new column A
$ws.Range("A1").Value = "Prefixo Cód. SCN"

# Fill in the numeric "Prefixo Cód. SCN" codes for rows 2-24
$codes = @(10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31801,31802)
for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $codes[$i]
}

# Update selection to match the target view
$ws.Range("B15").Select() | Out-Null
